$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARC_CO BM Prices")

# Update the publication date string in A3
$ws.Range("A3").Value = "April 10, 2025 1/"

# Update benchmark/actual price values in columns L and O for the affected rows
$updates = @{
    9  = 5.5
    10 = 6.5
    11 = 3.45
    12 = 0.255
    13 = 4.3499999999999996
    14 = 4.0999999999999996
    15 = 9.9499999999999993
    17 = 0.34799999999999998
    18 = 0.19900000000000001
    21 = 0.2155
    22 = 12.3
    23 = 0.48249999999999998
    24 = 0.2
    25 = 0.30299999999999999
    26 = 0.24
    28 = 0.33610000000000001
    29 = 0.14199999999999999
    30 = 0.152
    31 = 0.22500000000000001
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $ws.Range("L$row").Value = $value
    $ws.Range("O$row").Value = $value
}
